# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's existing
# "45.578.05"-style thousand-dotted text cells) instead of auto-coercing
# them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.514.95'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '2.379.28'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''333.11'
$ws.Range('E5').Value = '  +4.14%  '
$ws.Range('D6').Value = '''109.30'
$ws.Range('E6').Value = '  -6.01%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '''0.618'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').Value = '''41.16'
$ws.Range('E10').Value = '  -4.90%  '
$ws.Range('D11').Value = '''0.0922'
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('E12').Value = '  -2.60%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('E14').Value = '  -4.08%  '
$ws.Range('D15').Value = '2.743.60'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '''15.52'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '2.381.84'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '45.492.11'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '''15.25'
$ws.Range('E19').Value = '  +13.03%  '
$ws.Range('D20').Value = '''7.34'
$ws.Range('E20').Value = '  -3.82%  '
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('D22').Value = '''3.70'
$ws.Range('E22').Value = '  +3.20%  '
$ws.Range('D23').Value = '''73.44'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('D24').Value = '''264.96'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('D25').Value = '''2.33'
$ws.Range('E25').Value = '  -2.84%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').Value = '''7.50'
$ws.Range('E28').Value = '  -2.20%  '
$ws.Range('E29').Value = '  -1.92%  '
$ws.Range('D30').Value = '''22.54'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').Value = '''0.0958'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').Value = '''37.52'
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('D33').Value = '''169.51'
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('D35').Value = '''3.31'
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('D39').Value = '''1.97'
$ws.Range('E39').Value = '  +8.86%  '
$ws.Range('D40').Value = '''4.04'
$ws.Range('E40').Value = '  -5.48%  '
$ws.Range('D41').Value = '''0.0356'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('D42').Value = '''98.67'
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('D43').Value = '''70.95'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('E44').Value = '  -4.36%  '
$ws.Range('E45').Value = '  -4.12%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '''6.10'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.872.37'
$ws.Range('E47').Value = '  +12.93%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').Value = '''85.96'
$ws.Range('E49').Value = '  +4.24%  '
$ws.Range('D50').Value = '''9.40'
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('D51').Value = '''112.84'
$ws.Range('E51').Value = '  -4.19%  '
